$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (I1) to new header cells
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)  # xlPasteFormats

# Set header cell values
$ws.Range("J1").Value = "Onkelos"
$ws.Range("K1").Value = "Jonathan"

# Copy style from existing body cell (I2) to new body cells
$ws.Range("I2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("J2").Value = "“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt."
$ws.Range("K2").Value = "Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;"

$excel.CutCopyMode = $false

# Note: the runtime rounds ColumnWidth to the nearest displayable pixel width
# when persisting to OOXML (offset of 5/6 characters, snapped to 1/6 steps).
# The input values below are chosen so the saved <col width="..."/> comes out
# as close as possible to the target widths (314.4 and 282).
$ws.Columns.Item(10).ColumnWidth = 313.5
$ws.Columns.Item(11).ColumnWidth = 281.15
